# edit.ps1 - apply the commit's changes to the presentation
#
# Summary of changes (see diff):
#  - Slide 14: reposition/resize the DNS-history screenshot picture.
#  - Slide 23: reposition/resize the root-servers picture (tiny nudge) and
#              enlarge+grow the "https://root-servers.org/" citation textbox,
#              bumping its font size from 8pt to 14pt.
#  - Slide 25 ("Announcements"): tweak "9AM" -> "9AM and 9pm" and add a new
#              "Quiz 4 will be released today (due within 48 hours)." bullet.
#  - Slide 8: reposition/resize the Akamai citation textbox (font 8pt -> 14pt)
#              and reposition/resize the Akamai map picture.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 14 - "DNS: History" - resize/reposition the screenshot picture
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$pic14 = $s14.Shapes.Item(3)
$pic14.Left   = 393.0
$pic14.Top    = 0.34385826771653544
$pic14.Width  = 306.0
$pic14.Height = 530.4584251968504

# ---------------------------------------------------------------------------
# Slide 23 - "13 DNS root servers"
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Item(23)

# Picture: tiny nudge/resize
$pic23 = $s23.Shapes.Item(2)
$pic23.Left   = 54.00007874015748
$pic23.Top    = 141.47937007874015
$pic23.Width  = 623.5251181102362
$pic23.Height = 316.80001

# Citation textbox: grow the box and bump the font size to 14pt
$txt23 = $s23.Shapes.Item(6)
$txt23.Width  = 174.20985251968503
$txt23.Height = 24.234419448818898
$txt23.TextFrame.TextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# Slide 25 - "Announcements"
# ---------------------------------------------------------------------------
$s25 = $p.Slides.Item(25)
$body25 = $s25.Shapes.Item(6)
$tr25 = $body25.TextFrame.TextRange

$para9am = $tr25.Paragraphs(2)
$chars9am = $tr25.Characters($para9am.Start, $para9am.Length)
$chars9am.Text = "9AM and 9pm"

$paraAssignment = $tr25.Paragraphs(3)
$paraAssignment.InsertAfter("`r" + "Quiz 4 will be released today (due within 48 hours).")

# ---------------------------------------------------------------------------
# Slide 8 - "CDN example - Akamai"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# Citation textbox: reposition, widen and bump the font size to 14pt
$txt8 = $s8.Shapes.Item(5)
$txt8.Left   = 13.267175354330709
$txt8.Top    = 467.7655905511811
$txt8.Width  = 682.7328346456693
$txt8.Height = 24.234419448818898
$txt8.TextFrame.TextRange.Font.Size = 14

# Picture: reposition/resize
$pic8 = $s8.Shapes.Item(6)
$pic8.Left   = 24.0
$pic8.Top    = 119.55251968503937
$pic8.Width  = 654.0
$pic8.Height = 345.0663092125984
